# "mise a jour excel"
# Update the "Fonctionalités" tracking sheet:
#  - Row 10 (Admin / "Ordre des questions et réponses") is now done: Fait (O/N) -> O,
#    and the "Manque" note ("Logique à implémenter") is cleared.
#  - Row 16 (Admin / "Edition d'un questionnaire=>question") is now done: Fait (O/N) -> O,
#    and the "Manque" note ("Ordre") is cleared.
#  - Highlight the now-finished "Edition d'un questionnaire" / "Gestion des réponses" rows
#    (A16:F17) with the same green used elsewhere in the sheet to flag completed work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fonctionalités")
$ws.Activate()

# Row 10: "Ordre des questions et réponses" is done.
$ws.Range("D10").Value = "O"
$ws.Range("E10").Value = ""

# Row 16: "Edition d'un questionnaire=>question" is done.
$ws.Range("D16").Value = "O"
$ws.Range("E16").Value = ""

# Highlight rows 16-17 (Edition d'un questionnaire / Gestion des réponses) in green,
# matching the "Accent 6" theme color already used for the done-row conditional format.
$ws.Range("A16:F17").Interior.Color = 4697456

# Restore the user's on-screen selection at save time.
$ws.Range("B27").Select()
